# Update Name of Algo
# Apply updated RandomForest imputation results to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value  = 11.7831
$ws.Range("C8").Value  = -10.93649999999999
$ws.Range("C10").Value = -13.00700000000001
$ws.Range("C12").Value = -14.41960000000001
$ws.Range("D13").Value = -7.785200000000003
$ws.Range("C18").Value = -14.30670000000001
$ws.Range("E20").Value = 13.14749999999999
